$wb = $excel.ActiveWorkbook

# Remove "Chicken Nuggets" (row 4) and "Frozen Pizza" (row 2) from the "Quick" sheet.
# Delete the lower row first so the earlier row index stays valid.
$wsQuick = $wb.Worksheets.Item("Quick")
$wsQuick.Rows.Item(4).Delete()
$wsQuick.Rows.Item(2).Delete()

# Update the selection on the "Quick" sheet.
$wsQuick.Range("A10").Select()

# Make "Pastas" the active sheet/tab (was "Sides").
$wsPastas = $wb.Worksheets.Item("Pastas")
$wsPastas.Activate()
